$wb = $excel.ActiveWorkbook

# Add the new "prefaultvoltages" sheet after the last existing sheet ("linesequ")
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "prefaultvoltages"

# Header row
$ws.Range("A1").Value = "Bus"
$ws.Range("B1").Value = "PhaseA(Re)"
$ws.Range("C1").Value = "PhaseA(Im)"
$ws.Range("D1").Value = "PhaseB(Re)"
$ws.Range("E1").Value = "PhaseB(Im)"
$ws.Range("F1").Value = "PhaseC(Re)"
$ws.Range("G1").Value = "PhaseC(Im)"

# Make the new sheet the active / selected tab, matching the saved selection
$ws.Activate() | Out-Null
$ws.Range("J8").Select() | Out-Null
